# Applies the "retrait style dans html" edit to the SEO audit table:
# replaces rows 5-15 (categories B-F) with the refreshed set of
# SEO/accessibility remarks, widens columns B/C slightly, and moves
# the active selection to E15 (scrolled so C7 is the top-left cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "couleur texte pas assez visible"
$ws.Range("C5").Value = "changer la taille du texte"
$ws.Range("D5").Value = "on ne voit pas bien ce qu’il y a écrit "
$ws.Range("E5").Value = "changer la taille du texte"
$ws.Range("F5").Value = "expérience personnelle"

$ws.Range("B6").Value = "« accueil » non visible sur « contact »"
$ws.Range("C6").Value = "mettre un margin-right"
$ws.Range("D6").Value = "on ne voit pas le mot accueil en entier"
$ws.Range("E6").Value = "mettre un margin-right"

$ws.Range("B7").Value = "trop d’annuaires "
$ws.Range("C7").Value = "enlever des annuaires "
$ws.Range("D7").Value = "il y a trop d’annuaires sur la page"
$ws.Range("E7").Value = "mettre un minimum d’annuaires"

$ws.Range("B8").Value = "formulaire à refaire sur « contact »"
$ws.Range("C8").Value = "agrandir le formulaire"
$ws.Range("D8").Value = "le formulaire n’est pas assez lisible"
$ws.Range("E8").Value = "agrandir les formulaires"

$ws.Range("B9").Value = "« page2 » sur page contact"
$ws.Range("C9").Value = "enlever « page2 »"
$ws.Range("D9").Value = "on revient sur la même page"
$ws.Range("E9").Value = "retirer « page2 »"

$ws.Range("B10").Value = "point sur liste contact"
$ws.Range("C10").Value = "enlever les points de décoration"
$ws.Range("E10").Value = "mettre un list-style-decoration:none"

$ws.Range("B11").Value = "texte caché"
$ws.Range("C11").Value = "changer la couleur"
$ws.Range("D11").Value = "on ne voit pas le texte"
$ws.Range("E11").Value = "changer la couleur du texte"

$ws.Range("B12").Value = "Accueil sur index html"
$ws.Range("C12").Value = "enlever Accueil"
$ws.Range("D12").Value = "on s’y trouve déjà"
$ws.Range("E12").Value = "enlever Accueil"

$ws.Range("B13").Value = "le point sur H1"
$ws.Range("C13").Value = "enlever le point à la fin du titre"
$ws.Range("E13").Value = "enlever le point "

$ws.Range("B14").Value = "Liste 1 et 2"
$ws.Range("C14").Value = "enlever liste 1 et 2 sur annuaires"
$ws.Range("D14").Value = "on n’en a pas besoin"
$ws.Range("E14").Value = "retirer liste 1 et 2"

$ws.Range("B15").Value = "couleur icône réseaux"
$ws.Range("C15").Value = "changer la couleur des icônes "
$ws.Range("D15").Value = "cela ne se voit pas très bien en blanc"
$ws.Range("E15").Value = "changer la couleur"

# Column width tweaks (B: 30.94 -> 31.52, C: 34.65 -> 34.66 in the source file).
$ws.Columns.Item(2).ColumnWidth = 31.52
$ws.Columns.Item(3).ColumnWidth = 34.66

# Restore the cursor/scroll position recorded in the saved view.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 3
$ws.Range("E15").Select()
